$wb = $excel.ActiveWorkbook

# --- Update status text "Ready for handoff" -> "In Translation" everywhere it appears ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Narrow the Status-related columns (zh-cn/de-de width, and the corresponding
#     zh-cn/de-de status columns on the Overview sheet) to reflect the shorter text ---
$newWidth = 13.4101845877511 - (5 / 6)

$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth

$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
